$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 265, shifting existing rows 265:290 down to 267:292
$ws.Range("A265:A266").EntireRow.Insert()

# Fill in the first new row (265) - weekly price record added for week of 2022-01-17
$ws.Range("A265").Value = 10
$ws.Range("B265").Value = "Vega Modelo de Temuco"
$ws.Range("C265").Value = "La Araucanía"
$ws.Range("D265").Value = 44578
$ws.Range("E265").Value = 9
$ws.Range("F265").Value = 100112037
$ws.Range("G265").Value = "Cebollín"
$ws.Range("H265").Value = "Sin especificar"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 50
$ws.Range("K265").Value = 7000
$ws.Range("L265").Value = 8000
$ws.Range("M265").Value = 7600
$ws.Range("N265").Value = "$/docena de paquetes"
$ws.Range("O265").Value = "Provincia de Cautín"
$ws.Range("P265").Value = 633
$ws.Range("Q265").Value = 12
$ws.Range("R265").Value = "Hortaliza"

# Fill in the second new row (266)
$ws.Range("A266").Value = 10
$ws.Range("B266").Value = "Vega Modelo de Temuco"
$ws.Range("C266").Value = "La Araucanía"
$ws.Range("D266").Value = 44578
$ws.Range("E266").Value = 9
$ws.Range("F266").Value = 100112037
$ws.Range("G266").Value = "Cebollín"
$ws.Range("H266").Value = "Sin especificar"
$ws.Range("I266").Value = "Primera"
$ws.Range("J266").Value = 80
$ws.Range("K266").Value = 5000
$ws.Range("L266").Value = 5000
$ws.Range("M266").Value = 5000
$ws.Range("N266").Value = "$/docena de paquetes"
$ws.Range("O266").Value = "Región de O'Higgins"
$ws.Range("P266").Value = 417
$ws.Range("Q266").Value = 12
$ws.Range("R266").Value = "Hortaliza"
